# Add a new person ("Horst") to the weighted adjacency matrix and update
# several edge weights from 1 -> 2 (upgrade) or 0 -> 2 (new weighted edge).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Matrix")

# --- Add the new person "Horst" as row 11 / column K ---
$ws.Range("A11").Value = "Horst"
$ws.Range("K1").Value = "Horst"

# --- Clear the (self-loop) diagonal zero cells that sat directly right of each name ---
$ws.Range("B2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("F6").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("J10").ClearContents()

# --- Update existing edge weights ---
$ws.Range("E2").Value = 2
$ws.Range("I2").Value = 2
$ws.Range("D3").Value = 2
$ws.Range("I3").Value = 2
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 2
$ws.Range("H5").Value = 2
$ws.Range("I6").Value = 2
$ws.Range("J8").Value = 2

# --- New column K (weights against "Horst") ---
$ws.Range("K2").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("K4").Value = 2
$ws.Range("K5").Value = 1
$ws.Range("K6").Value = 2
$ws.Range("K7").Value = 0
$ws.Range("K8").Value = 1
$ws.Range("K9").Value = 2
$ws.Range("K10").Value = 0

# --- View / selection cosmetics matching the authored file ---
$ws.Range("K11").Select() | Out-Null
